$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the time_taken timestamps on the "data" sheet (F2:F14) ---
$newTimes = @(
  "2021-10-05 14:22:11.112192",
  "2021-10-05 14:22:11.112203",
  "2021-10-05 14:22:11.112208",
  "2021-10-05 14:22:11.112212",
  "2021-10-05 14:22:11.112215",
  "2021-10-05 14:22:11.112218",
  "2021-10-05 14:22:11.112221",
  "2021-10-05 14:22:11.112223",
  "2021-10-05 14:22:11.112226",
  "2021-10-05 14:22:11.112229",
  "2021-10-05 14:22:11.112232",
  "2021-10-05 14:22:11.112235",
  "2021-10-05 14:22:11.112237"
)
for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add a new "metadata" worksheet positioned after the "data" sheet ---
$ws = $wb.Worksheets.Add([System.Type]::Missing, $dataSheet)
$ws.Name = "metadata"

# Header row (B1:G1) - match header style used on the "data" sheet
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row (A2:G2)
$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Peeling skin syndrome"
$ws.Range("C2").Value = 24
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.2"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "2017-11-05T02:37:19.863085Z"
$ws.Range("F2").Value = "2021-10-05 14:22:11.108785"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/24/?format=json"
